$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidential note date (2021-04-08 -> 2021-04-09), preserving the line break
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values
$ws.Range("D2").Value = 0.2530612750269134
$ws.Range("E2").Value = 0.006614850339011102

$ws.Range("D3").Value = 0.4904784278233771
$ws.Range("E3").Value = 0.006548431105047658

$ws.Range("D4").Value = 0.1011918147985575
$ws.Range("E4").Value = 0.00746268656716409

$ws.Range("D5").Value = 0.09837409844041672
$ws.Range("E5").Value = 0.005160550458715552

$ws.Range("D6").Value = 0.0568943839107352
$ws.Range("E6").Value = 0.001376462491397001

$ws.Range("E7").Value = 0.006226966935276801
